# Modified Generic Class functions and included Runmodes.
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "TestCases" worksheet right after "TestData" (i.e.
#    before "ListConfig") and populate it with the TestNames/RunMode
#    table plus a Y/N list-validation sourced from ListConfig!D1:D2.
# ------------------------------------------------------------------
$testData = $wb.Worksheets.Item("TestData")
$ws = $wb.Worksheets.Add($null, $testData)
$ws.Name = "TestCases"

$ws.Range("A1").Value = "TestNames"
$ws.Range("B1").Value = "RunMode"
$ws.Range("A2").Value = "LoginAsBankManager"
$ws.Range("B2").Value = "Y"
$ws.Range("A3").Value = "AddCustomers"
$ws.Range("B3").Value = "Y"
$ws.Range("A4").Value = "OpenAccount"
$ws.Range("B4").Value = "N"

$ws.Columns.Item(1).ColumnWidth = 20.140625
$ws.Columns.Item(2).ColumnWidth = 11.140625

$ws.Range("A1:B1").Style = "Header"
$ws.Range("A2:A4").Style = "Bordered"
$ws.Range("B2:B4").Style = "BorderedCentered"

$ws.Range("B2:B4").Validation.Add(3, 1, 1, '=ListConfig!$D$1:$D$2')

# ------------------------------------------------------------------
# 2. Add a RunMode (Y/N) helper column to "ListConfig" so the new
#    validation above has something to point at.
# ------------------------------------------------------------------
$listConfig = $wb.Worksheets.Item("ListConfig")
$listConfig.Range("D1").Value = "Y"
$listConfig.Range("D2").Value = "N"
$listConfig.Range("D1:D2").Style = "Header"
$listConfig.Range("D1").Select()

# ------------------------------------------------------------------
# 3. Make "TestCases" the active/selected sheet (matches activeTab).
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1").Select()
